# Add the new "ORDEN DE DOMICILIACIÓN DE ADEUDO DIRECTO SEPA" section
# (rows 41-48) to the bottom of the "FICHA CLIENTE" sheet, mirroring the
# look & feel of the existing sections above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Cell values (order matters: it drives shared-string allocation) ---
$ws.Range("A42").Value = "ORDEN DE DOMICILIACIÓN DE ADEUDO DIRECTO SEPA"
$ws.Range("A43").Value = "Nombre entidad bancaria"
$ws.Range("A44").Value = "Domicilio entidad bancaria"
$ws.Range("A45").Value = "Código postal"
$ws.Range("A46").Value = "Población"
$ws.Range("A47").Value = "Provincia"
$ws.Range("A48").Value = "Nº DE CUENTA (IBAN 24 caracteres)"

# --- 2. Formatting: reuse the look of the existing analogous rows ---
# Row 42 is a section-title bar, like rows 11 / 16 / 37 ("CONTACTOS", etc.)
$ws.Range("A11:D11").Copy()
$ws.Range("A42:D42").PasteSpecial(-4122)

# Rows 43-48 are label/input pairs, like rows 18 / 22 / 27 / 32 ("Nombre")
$ws.Range("A18:B18").Copy()
$ws.Range("A43:B43").PasteSpecial(-4122)
$ws.Range("A44:B44").PasteSpecial(-4122)
$ws.Range("A45:B45").PasteSpecial(-4122)
$ws.Range("A46:B46").PasteSpecial(-4122)
$ws.Range("A47:B47").PasteSpecial(-4122)
$ws.Range("A48:B48").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Row heights to match the rest of the sheet (15pt, thick-border rows)
$ws.Rows("41:48").RowHeight = 15

# --- 3. Column A is now much wider to fit the new longer labels ---
$ws.Columns("A").ColumnWidth = 45.109375

# --- 4. View state: scrolled down to the new section, B43 selected ---
$ws.Application.ActiveWindow.ScrollRow = 32
$ws.Range("B43").Select()

# --- 5. Basic page setup for printing ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
